# 18/02/2026 - Reportes Actualizados, Acta de Validacion PDF y Word actualizados,
# Acta de Computo Total PDF y Word agregados
#
# The sheet "Hoja1" gains a new blank row at the very top (pushing every
# existing row down by one) and loses two rows that were already blank
# spacer rows (the old row 5, which only had an empty F5 cell, and the old
# row 7 -- which, after the insert above, sits at row 8 -- an empty F-only
# spacer row). Net effect on row numbering, from the bottom of the sheet
# upward: row 12 -> 11, row 11 -> 10, ... row 6 -> 6 (unchanged), row 4 -> 5,
# row 3 -> 4, row 2 -> 3, row 1 -> 2, and a brand new empty row 1 appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new blank row above the current row 1 - shifts everything down.
$ws.Rows("1").Insert()

# Remove the (now shifted) empty spacer rows. Delete the lower-numbered
# row first so the second row's index doesn't move before we get to it.
$ws.Rows("8").Delete()
$ws.Rows("6").Delete()
